$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72..129 down to 73..130
$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new record's data.
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
#  Calidad, Origen, Clasificacion mirror the row that used to sit at 72.)
$ws.Range("A72").Value = 5
$ws.Range("B72").Value = "Macroferia Regional de Talca"
$ws.Range("C72").Value = "Maule"
$ws.Range("D72").Value = 44512
$ws.Range("E72").Value = 7
$ws.Range("F72").Value = 100112021
$ws.Range("G72").Value = "Ají"
$ws.Range("H72").Value = "Americana (o)"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 150
$ws.Range("K72").Value = 23000
$ws.Range("L72").Value = 23000
$ws.Range("M72").Value = 23000
$ws.Range("N72").Value = "$/caja 15 kilos"
$ws.Range("O72").Value = "Región del Maule"
$ws.Range("P72").Value = 1533
$ws.Range("Q72").Value = 15
$ws.Range("R72").Value = "Hortaliza"
